$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-6 hold daily price records for "Caqui" (persimmon) at the same
# market; the edit re-shuffles which date/variety/quality/price/origin block
# belongs to which row (a weekly re-sequencing of the same underlying data),
# per the commit message "Fruta / hortaliza, semanal".
#
# Snapshot the full row (columns A..T) for each of rows 2..6 first, then
# write the rows back out in rotated order so every column lands on the
# value the diff shows, without us needing to special-case each column.

$lastCol = 20  # column T

function Get-RowSnapshot($rowIndex) {
    $vals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $vals += , ($ws.Cells.Item($rowIndex, $c).Value())
    }
    return $vals
}

function Set-RowSnapshot($rowIndex, $vals) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($rowIndex, $c).Value = $vals[$c - 1]
    }
}

$row2 = Get-RowSnapshot 2
$row3 = Get-RowSnapshot 3
$row4 = Get-RowSnapshot 4
$row5 = Get-RowSnapshot 5
$row6 = Get-RowSnapshot 6

# Cycle (2 4 5): row2 <- old row4, row4 <- old row5, row5 <- old row2
Set-RowSnapshot 2 $row4
Set-RowSnapshot 4 $row5
Set-RowSnapshot 5 $row2

# Cycle (3 6): row3 <- old row6, row6 <- old row3
Set-RowSnapshot 3 $row6
Set-RowSnapshot 6 $row3
